# Apply numeric corrections to the profit-tracking sheets.
# Values were recomputed upstream by the scheduled runner; this script
# writes the refreshed figures (and clears/adds a few cells whose
# presence changed) cell-by-cell per sheet.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1239.2727
$ws.Range("J17").Value = 1247.6666
$ws.Range("L17").Value = 3742.9998
$ws.Range("N17").Value = -4078.9998
$ws.Range("H42").Value = 96.5
$ws.Range("I42").Value = 73.53846
$ws.Range("J42").Value = 196
$ws.Range("K42").Value = 220.61538
$ws.Range("L42").Value = 588
$ws.Range("M42").Value = 9.384619999999984
$ws.Range("N42").Value = -1048
$ws.Range("H64").Value = 1950
$ws.Range("H67").Value = 1950
$ws.Range("H76").Value = 4999.5
$ws.Range("I76").Value = 4999.5
$ws.Range("K76").Value = 4999.5
$ws.Range("M76").Value = -4684.5
$ws.Range("H79").Value = 4999.5
$ws.Range("I79").Value = 4999.5
$ws.Range("K79").Value = 4999.5
$ws.Range("M79").Value = -3907.5
$ws.Range("H100").Value = 1654.8572
$ws.Range("I100").Value = 1433.8182
$ws.Range("J100").Value = 2465.3333
$ws.Range("K100").Value = 1433.8182
$ws.Range("L100").Value = 2465.3333
$ws.Range("M100").Value = -892.8181999999999
$ws.Range("N100").Value = -3547.3333
$ws.Range("H116").Value = 2000
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 1755.75
$ws.Range("I132").Value = 1720.8572
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5162.571599999999
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2632.571599999999
$ws.Range("N132").Value = -11060
$ws.Range("H137").Value = 1625.3572
$ws.Range("I137").Value = 1475.6
$ws.Range("K137").Value = 4426.799999999999
$ws.Range("M137").Value = -1876.799999999999
$ws.Range("H141").Value = 1978.8948
$ws.Range("I141").Value = 1255.5278
$ws.Range("K141").Value = 3766.5834
$ws.Range("M141").Value = 1413.4166

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8453.058000000001
$ws.Range("I32").Value = 6071.1377
$ws.Range("K32").Value = 6071.1377
$ws.Range("M32").Value = -5784.1377
$ws.Range("H45").Value = 2503.5386
$ws.Range("I45").Value = 2503.8333
$ws.Range("K45").Value = 2503.8333
$ws.Range("M45").Value = -2126.8333
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H74").Value = 1150
$ws.Range("I74").Value = 1150
$ws.Range("K74").Value = 1150
$ws.Range("M74").Value = -276
$ws.Range("H76").Value = 16457.8
$ws.Range("J76").Value = 16457.8
$ws.Range("L76").Value = 16457.8
$ws.Range("N76").Value = -17133.8
$ws.Range("H77").Value = 1150
$ws.Range("I77").Value = 1150
$ws.Range("K77").Value = 5750
$ws.Range("M77").Value = -1382
$ws.Range("H79").Value = 16457.8
$ws.Range("J79").Value = 16457.8
$ws.Range("L79").Value = 16457.8
$ws.Range("N79").Value = -18797.8
$ws.Range("H88").Value = 1180
$ws.Range("J88").Value = 1600
$ws.Range("L88").Value = 1600
$ws.Range("N88").Value = -2412
$ws.Range("H91").Value = 1180
$ws.Range("J91").Value = 1600
$ws.Range("L91").Value = 1600
$ws.Range("N91").Value = -4408
$ws.Range("H97").Value = 1090.125
$ws.Range("I97").Value = 1074.4286
$ws.Range("K97").Value = 1074.4286
$ws.Range("M97").Value = -578.4286
$ws.Range("H132").Value = 1325.2354
$ws.Range("I132").Value = 1324.2142
$ws.Range("J132").Value = 1330
$ws.Range("K132").Value = 3972.6426
$ws.Range("L132").Value = 3990
$ws.Range("M132").Value = -1442.6426
$ws.Range("N132").Value = -9050
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1581
$ws.Range("I86").Value = 1349.5
$ws.Range("K86").Value = 1349.5
$ws.Range("M86").Value = -226.5
$ws.Range("H89").Value = 1581
$ws.Range("I89").Value = 1349.5
$ws.Range("K89").Value = 6747.5
$ws.Range("M89").Value = -1131.5
$ws.Range("H105").Value = 2216.5
$ws.Range("I105").Value = 2116.375
$ws.Range("J105").Value = 2617
$ws.Range("K105").Value = 2116.375
$ws.Range("L105").Value = 2617
$ws.Range("M105").Value = -369.375
$ws.Range("N105").Value = -6111
$ws.Range("H134").Value = 2867.6365
$ws.Range("J134").Value = 3981.818
$ws.Range("L134").Value = 11945.454
$ws.Range("N134").Value = -17015.454

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 250
$ws.Range("K22").Value = 250
$ws.Range("M22").Value = 100
$ws.Range("H58").Value = 2026.7742
$ws.Range("I58").Value = 1146.125
$ws.Range("J58").Value = 5046.143
$ws.Range("K58").Value = 1146.125
$ws.Range("L58").Value = 5046.143
$ws.Range("M58").Value = -943.125
$ws.Range("N58").Value = -5452.143
$ws.Range("H62").Value = 102476
$ws.Range("I62").Value = 3301.6667
$ws.Range("K62").Value = 3301.6667
$ws.Range("M62").Value = -2677.6667
$ws.Range("H65").Value = 102476
$ws.Range("I65").Value = 3301.6667
$ws.Range("K65").Value = 16508.3335
$ws.Range("M65").Value = -13388.3335
$ws.Range("H80").Value = 29999.5
$ws.Range("I80").Value = 29999
$ws.Range("K80").Value = 29999
$ws.Range("M80").Value = -28876
$ws.Range("H83").Value = 29999.5
$ws.Range("I83").Value = 29999
$ws.Range("K83").Value = 89997
$ws.Range("M83").Value = -84381
$ws.Range("H86").Value = 10128.429
$ws.Range("I86").Value = 8974.75
$ws.Range("J86").Value = 11666.667
$ws.Range("K86").Value = 8974.75
$ws.Range("L86").Value = 11666.667
$ws.Range("M86").Value = -7851.75
$ws.Range("N86").Value = -13912.667
$ws.Range("H89").Value = 10128.429
$ws.Range("I89").Value = 8974.75
$ws.Range("J89").Value = 11666.667
$ws.Range("K89").Value = 44873.75
$ws.Range("L89").Value = 58333.335
$ws.Range("M89").Value = -39257.75
$ws.Range("N89").Value = -69565.33499999999
$ws.Range("H122").Value = 2563.7058
$ws.Range("I122").Value = 2563.7058
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7691.117400000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5241.117400000001
$ws.Range("N122").ClearContents()
$ws.Range("H130").Value = 54250
$ws.Range("I130").Value = 49000
$ws.Range("K130").Value = 49000
$ws.Range("M130").Value = -43980
$ws.Range("H136").Value = 2026.7742
$ws.Range("I136").Value = 1146.125
$ws.Range("J136").Value = 5046.143
$ws.Range("K136").Value = 3438.375
$ws.Range("L136").Value = 15138.429
$ws.Range("M136").Value = -888.375
$ws.Range("N136").Value = -20238.429

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 25000
$ws.Range("J15").Value = 25000
$ws.Range("L15").Value = 25000
$ws.Range("N15").Value = -25576
$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996
$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984
$ws.Range("H101").Value = 84292
$ws.Range("J101").Value = 84292
$ws.Range("L101").Value = 84292
$ws.Range("N101").Value = -90782

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2224
$ws.Range("I7").Value = 2224
$ws.Range("K7").Value = 2224
$ws.Range("M7").Value = -2112
$ws.Range("H61").Value = 4335.3076
$ws.Range("I61").Value = 4305.364
$ws.Range("K61").Value = 4305.364
$ws.Range("M61").Value = -4103.364
$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 1000
$ws.Range("M100").Value = -459
$ws.Range("H113").Value = 4335.3076
$ws.Range("I113").Value = 4305.364
$ws.Range("K113").Value = 4305.364
$ws.Range("M113").Value = -2135.364
$ws.Range("H126").Value = 2224
$ws.Range("I126").Value = 2224
$ws.Range("K126").Value = 6672
$ws.Range("M126").Value = -4202
$ws.Range("H132").Value = 2242.9143
$ws.Range("I132").Value = 2016.3226
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 6048.9678
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -3518.9678
$ws.Range("N132").Value = -17057

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2125.125
$ws.Range("I122").Value = 2699.75
$ws.Range("J122").Value = 1933.5834
$ws.Range("K122").Value = 8099.25
$ws.Range("L122").Value = 5800.7502
$ws.Range("M122").Value = -5649.25
$ws.Range("N122").Value = -10700.7502
$ws.Range("H136").Value = 1665.4117
$ws.Range("I136").Value = 901
$ws.Range("K136").Value = 2703
$ws.Range("M136").Value = -153

